$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 22.76743878352112
$ws.Range("D2").Value = 9.140147260761253
$ws.Range("E2").Value = 11.3745074639273
$ws.Range("F2").Value = 58.29784711056954
$ws.Range("G2").Value = 3.782636511991054
$ws.Range("I2").Value = 30.65602340278918
$ws.Range("J2").Value = 7.924709575475915
$ws.Range("L2").Value = 19.79125620569425
$ws.Range("N2").Value = 19.66492604090251
$ws.Range("B3").Value = 22.46943343895988
$ws.Range("D3").Value = 8.808267140785949
$ws.Range("E3").Value = 10.82820605269461
$ws.Range("F3").Value = 58.00964452658879
$ws.Range("G3").Value = 3.790019376043156
$ws.Range("I3").Value = 30.80469895541967
$ws.Range("J3").Value = 7.715299493542672
$ws.Range("L3").Value = 19.58109222547393
$ws.Range("N3").Value = 19.75111559920293
$ws.Range("B4").Value = 22.29033688655862
$ws.Range("D4").Value = 8.602129712792191
$ws.Range("E4").Value = 10.48261520367097
$ws.Range("F4").Value = 57.85272152179812
$ws.Range("G4").Value = 3.794774951946768
$ws.Range("I4").Value = 30.90319957693454
$ws.Range("J4").Value = 7.582686929866058
$ws.Range("L4").Value = 19.45674804697492
$ws.Range("N4").Value = 19.80615161317976
$ws.Range("B5").Value = 22.2184101678976
$ws.Range("D5").Value = 8.517679649608697
$ws.Range("E5").Value = 10.33944143732984
$ws.Range("F5").Value = 57.79382718427043
$ws.Range("G5").Value = 3.796769137401167
$ws.Range("I5").Value = 30.94514408392297
$ws.Range("J5").Value = 7.527669890929575
$ws.Range("L5").Value = 19.40730439047137
$ws.Range("N5").Value = 19.82911444107442
$ws.Range("B6").Value = 22.20653292974637
$ws.Range("D6").Value = 8.503634263276124
$ws.Range("E6").Value = 10.31553284813154
$ws.Range("F6").Value = 57.78435329402514
$ws.Range("G6").Value = 3.797103676730628
$ws.Range("I6").Value = 30.95221767504919
$ws.Range("J6").Value = 7.518476553226606
$ws.Range("L6").Value = 19.39916975549878
$ws.Range("N6").Value = 19.8329598226476
$ws.Range("B7").Value = 22.28936247420194
$ws.Range("D7").Value = 8.600992399193231
$ws.Range("E7").Value = 10.48069349500183
$ws.Range("F7").Value = 57.85190677860643
$ws.Range("G7").Value = 3.794801618066435
$ws.Range("I7").Value = 30.90375795985961
$ws.Range("J7").Value = 7.581948851655171
$ws.Range("L7").Value = 19.45607620268402
$ws.Range("N7").Value = 19.80645912650049
$ws.Range("B8").Value = 22.66393040967044
$ws.Range("D8").Value = 9.026289448363348
$ws.Range("E8").Value = 11.18836929015395
$ws.Range("F8").Value = 58.19432836332516
$ws.Range("G8").Value = 3.785136136889241
$ws.Range("I8").Value = 30.70578407451195
$ws.Range("J8").Value = 7.853364111025011
$ws.Range("L8").Value = 19.7178453874458
$ws.Range("N8").Value = 19.69420717641568
$ws.Range("B9").Value = 23.42560592882016
$ws.Range("D9").Value = 9.835726065192997
$ws.Range("E9").Value = 12.48766764655857
$ws.Range("F9").Value = 59.02372588726958
$ws.Range("G9").Value = 3.76793322454359
$ws.Range("I9").Value = 30.37521048139268
$ws.Range("J9").Value = 8.35208560854471
$ws.Range("L9").Value = 20.266285966311
$ws.Range("N9").Value = 19.49071235897537
$ws.Range("B10").Value = 23.99682981689741
$ws.Range("D10").Value = 10.40850428077383
$ws.Range("E10").Value = 13.46037264668452
$ws.Range("F10").Value = 59.7274016346478
$ws.Range("G10").Value = 3.75634190871166
$ws.Range("I10").Value = 30.16807566905796
$ws.Range("J10").Value = 8.696245916668701
$ws.Range("L10").Value = 20.68760239173575
$ws.Range("N10").Value = 19.35112603183493
$ws.Range("B11").Value = 24.25817148248617
$ws.Range("D11").Value = 10.66309620667537
$ws.Range("E11").Value = 13.92196127765591
$ws.Range("F11").Value = 60.06742573195956
$ws.Range("G11").Value = 3.751291908432276
$ws.Range("I11").Value = 30.08173802892676
$ws.Range("J11").Value = 8.847656425944466
$ws.Range("L11").Value = 20.88258827157543
$ws.Range("N11").Value = 19.28973155700286
$ws.Range("B12").Value = 24.35726268621883
$ws.Range("D12").Value = 10.75855928555937
$ws.Range("E12").Value = 14.09296522333941
$ws.Range("F12").Value = 60.19898841479765
$ws.Range("G12").Value = 3.749411323319951
$ws.Range("I12").Value = 30.05019024654166
$ws.Range("J12").Value = 8.904227420149256
$ws.Range("L12").Value = 20.95684529517258
$ws.Range("N12").Value = 19.26678190850678
$ws.Range("B13").Value = 24.33591735430159
$ws.Range("D13").Value = 10.73804292151959
$ws.Range("E13").Value = 14.05630468432616
$ws.Range("F13").Value = 60.17053041960207
$ws.Range("G13").Value = 3.749814934420368
$ws.Range("I13").Value = 30.05693344476434
$ws.Range("J13").Value = 8.892078222430246
$ws.Range("L13").Value = 20.94083496093869
$ws.Range("N13").Value = 19.27171127841849
$ws.Range("B14").Value = 24.26632172813874
$ws.Range("D14").Value = 10.67096941875818
$ws.Range("E14").Value = 13.93610570861174
$ws.Range("F14").Value = 60.07819367198424
$ws.Range("G14").Value = 3.75113655727751
$ws.Range("I14").Value = 30.0791195283846
$ws.Range("J14").Value = 8.852326015127179
$ws.Range("L14").Value = 20.88868935623336
$ws.Range("N14").Value = 19.28783750097855
$ws.Range("B15").Value = 24.22370634822908
$ws.Range("D15").Value = 10.62975944640485
$ws.Range("E15").Value = 13.8619874445516
$ws.Range("F15").Value = 60.02199784384713
$ws.Range("G15").Value = 3.751950213456476
$ws.Range("I15").Value = 30.09285880264774
$ws.Range("J15").Value = 8.827876352660242
$ws.Range("L15").Value = 20.85680164602861
$ws.Range("N15").Value = 19.29775412739628
$ws.Range("B16").Value = 23.97977298193168
$ws.Range("D16").Value = 10.39173830767096
$ws.Range("E16").Value = 13.42967533878691
$ws.Range("F16").Value = 59.70557597762767
$ws.Range("G16").Value = 3.756676397151963
$ws.Range("I16").Value = 30.17387780012279
$ws.Range("J16").Value = 8.686245115599879
$ws.Range("L16").Value = 20.67492157080584
$ws.Range("N16").Value = 19.35518034953966
$ws.Range("B17").Value = 23.8304464652882
$ws.Range("D17").Value = 10.24412600219975
$ws.Range("E17").Value = 13.15769457736941
$ws.Range("F17").Value = 59.51652146810602
$ws.Range("G17").Value = 3.759632633640122
$ws.Range("I17").Value = 30.22560931149668
$ws.Range("J17").Value = 8.598021794886353
$ws.Range("L17").Value = 20.56415471144976
$ws.Range("N17").Value = 19.39094592696155
$ws.Range("B18").Value = 23.74470309238988
$ws.Range("D18").Value = 10.15866608374703
$ws.Range("E18").Value = 12.99876373917148
$ws.Range("F18").Value = 59.40966124090406
$ws.Range("G18").Value = 3.761353986786844
$ws.Range("I18").Value = 30.25610572385631
$ws.Range("J18").Value = 8.546794639513234
$ws.Range("L18").Value = 20.5007612022151
$ws.Range("N18").Value = 19.41171559230942
$ws.Range("B19").Value = 23.71569950257191
$ws.Range("D19").Value = 10.12963795802678
$ws.Range("E19").Value = 12.94775661980598
$ws.Range("F19").Value = 59.3738046721429
$ws.Range("G19").Value = 3.761940424537301
$ws.Range("I19").Value = 30.26655833297211
$ws.Range("J19").Value = 8.529367760221797
$ws.Range("L19").Value = 20.47935340965356
$ws.Range("N19").Value = 19.41878199851853
$ws.Range("B20").Value = 23.84632812929714
$ws.Range("D20").Value = 10.25989792123875
$ws.Range("E20").Value = 13.18690544637252
$ws.Range("F20").Value = 59.536452605323
$ws.Range("G20").Value = 3.75931576567825
$ws.Range("I20").Value = 30.22002554311873
$ws.Range("J20").Value = 8.607463536423825
$ws.Range("L20").Value = 20.57591368180886
$ws.Range("N20").Value = 19.38711812686174
$ws.Range("B21").Value = 24.28676089639276
$ws.Range("D21").Value = 10.69069683016809
$ws.Range("E21").Value = 13.9715137490737
$ws.Range("F21").Value = 60.10523966232347
$ws.Range("G21").Value = 3.750747505737437
$ws.Range("I21").Value = 30.07257172744057
$ws.Range("J21").Value = 8.864023141701505
$ws.Range("L21").Value = 20.90399484773341
$ws.Range("N21").Value = 19.28309274850315
$ws.Range("B22").Value = 24.57530909024392
$ws.Range("D22").Value = 10.96670563107336
$ws.Range("E22").Value = 14.46222152673291
$ws.Range("F22").Value = 60.49327974420869
$ws.Range("G22").Value = 3.74533250988092
$ws.Range("I22").Value = 29.98289008273679
$ws.Range("J22").Value = 9.027231904090328
$ws.Range("L22").Value = 21.12083655606481
$ws.Range("N22").Value = 19.21684794675618
$ws.Range("B23").Value = 24.42126979465991
$ws.Range("D23").Value = 10.81992777095727
$ws.Range("E23").Value = 14.20233500550711
$ws.Range("F23").Value = 60.28470542563146
$ws.Range("G23").Value = 3.748205786750948
$ws.Range("I23").Value = 30.03013886975645
$ws.Range("J23").Value = 8.940540489171863
$ws.Range("L23").Value = 21.00490172078964
$ws.Range("N23").Value = 19.25204581906899
$ws.Range("B24").Value = 23.8391476908463
$ws.Range("D24").Value = 10.25276928518558
$ws.Range("E24").Value = 13.17370720485802
$ws.Range("F24").Value = 59.52743603441598
$ws.Range("G24").Value = 3.759458953840034
$ws.Range("I24").Value = 30.22254761318509
$ws.Range("J24").Value = 8.603196500499639
$ws.Range("L24").Value = 20.57059655035036
$ws.Range("N24").Value = 19.38884802852787
$ws.Range("B25").Value = 23.21715401339163
$ws.Range("D25").Value = 9.620146249723922
$ws.Range("E25").Value = 12.14662913021642
$ws.Range("F25").Value = 58.78261622329368
$ws.Range("G25").Value = 3.772401647360075
$ws.Range("I25").Value = 30.4584074186308
$ws.Range("J25").Value = 8.220955634814169
$ws.Range("L25").Value = 20.11446559251021
$ws.Range("N25").Value = 19.54400501693766
